$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.174.77"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.915.94"
$ws.Range("E3").Value = "  +0.18%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.34%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "325.06"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  +0.13%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3849"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07795"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.47%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9678"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.27%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "22.27"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").Value = "1.918.51"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "7.039"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.735"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.53%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.07067"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "85.91"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("E17").Value = "  +0.19%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000009655"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.57%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "16.92"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "29.146.26"
$ws.Range("E21").Value = "  +0.28%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.456"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.59%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "11.04"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.154.17"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.089"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "157.28"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "19.28"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "5.697"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.62%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "118.03"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.817"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.09343"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.8549"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.114"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.270"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.072"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.05741"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.160"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02067"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.02%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "7.561"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5593"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.000003082"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +6.43%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1768"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "9.227"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.72%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.732"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +6.81%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.5236"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "11.38"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.06825"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.52%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.058"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -6.54%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.795"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "111.04"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.2984"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.97%  "
